# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Femacal de La Calera - Repollo"
# as row 324, shifting all subsequent rows (old 324-370) down by one
# (new 325-371).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 324; this pushes the existing data
# (previously rows 324:370) down to rows 325:371 and extends the
# sheet dimension from A1:R370 to A1:R371 automatically.
$ws.Rows("324:324").Insert()

# Populate the newly inserted row 324 with this week's record.
$ws.Range("A324").Value = 3
$ws.Range("B324").Value = "Femacal de La Calera"
$ws.Range("C324").Value = "Coquimbo"
$ws.Range("D324").Value = 44491
$ws.Range("E324").Value = 5
$ws.Range("F324").Value = 100112006
$ws.Range("G324").Value = "Repollo"
$ws.Range("H324").Value = "Crespo record"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 4650
$ws.Range("K324").Value = 500
$ws.Range("L324").Value = 650
$ws.Range("M324").Value = 585
$ws.Range("N324").Value = "$/unidad"
$ws.Range("O324").Value = "Provincia de Quillota"
$ws.Range("P324").Value = 585
$ws.Range("Q324").Value = 1
$ws.Range("R324").Value = "Hortaliza"
